# Edit script generated to match the target diff for decision_data_new.xlsx (sheet "sys2")
# - Re-points the _xlchart.v1.* defined names used by the box-whisker charts. Excel
#   renumbers these automatically on save; net effect on the charts (which physical
#   column each chart/series plots) is unchanged - only the internal index labels swap.
# - Rounds columns C and F in sys2 rows 88-101 to 9 decimal places.
# - Appends 17 new data rows (102-118) to the sys2 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sys2")

# --- Re-point the _xlchart.v1.* defined names (Excel renumbers these on save) ---
$wb.Names.Item("_xlchart.v1.1").RefersTo = "=wall_mounted_data!`$F`$1"
$wb.Names.Item("_xlchart.v1.10").RefersTo = "=wall_mounted_data!`$B`$1"
$wb.Names.Item("_xlchart.v1.11").RefersTo = "=wall_mounted_data!`$B`$2:`$B`$320"
$wb.Names.Item("_xlchart.v1.13").RefersTo = "=wall_mounted_data!`$E`$1"
$wb.Names.Item("_xlchart.v1.14").RefersTo = "=wall_mounted_data!`$E`$2:`$E`$320"
$wb.Names.Item("_xlchart.v1.2").RefersTo = "=wall_mounted_data!`$F`$2:`$F`$320"
$wb.Names.Item("_xlchart.v1.4").RefersTo = "=wall_mounted_data!`$C`$1"
$wb.Names.Item("_xlchart.v1.5").RefersTo = "=wall_mounted_data!`$C`$2:`$C`$320"
$wb.Names.Item("_xlchart.v1.7").RefersTo = "=wall_mounted_data!`$D`$1"
$wb.Names.Item("_xlchart.v1.8").RefersTo = "=wall_mounted_data!`$D`$2:`$D`$320"

# --- Re-round existing values in columns C and F (rows 88-101) to 9 decimals ---
$ws.Cells.Item(88, 3).Value = 0.08267205
$ws.Cells.Item(89, 3).Value = 0.045515433
$ws.Cells.Item(90, 3).Value = 0.041364132
$ws.Cells.Item(91, 3).Value = 0.051255658
$ws.Cells.Item(91, 6).Value = 0.054635584
$ws.Cells.Item(92, 3).Value = -0.066458714
$ws.Cells.Item(92, 6).Value = 0.05164132
$ws.Cells.Item(93, 3).Value = -0.10568711
$ws.Cells.Item(93, 6).Value = -0.016643454
$ws.Cells.Item(94, 3).Value = -0.086765499
$ws.Cells.Item(94, 6).Value = -0.016101219
$ws.Cells.Item(95, 3).Value = -0.077654472
$ws.Cells.Item(95, 6).Value = -0.011960002
$ws.Cells.Item(96, 3).Value = -0.111118108
$ws.Cells.Item(96, 6).Value = -0.011786788
$ws.Cells.Item(97, 3).Value = -0.062132519
$ws.Cells.Item(97, 6).Value = -0.163323181
$ws.Cells.Item(98, 3).Value = -0.080074349
$ws.Cells.Item(98, 6).Value = -0.024401327
$ws.Cells.Item(99, 3).Value = -0.085107872
$ws.Cells.Item(99, 6).Value = -0.011141382
$ws.Cells.Item(100, 3).Value = -0.086319378
$ws.Cells.Item(100, 6).Value = -0.035547655
$ws.Cells.Item(101, 3).Value = -0.068487588
$ws.Cells.Item(101, 6).Value = -0.021448924

# --- Append new rows 102-118 ---
$ws.Cells.Item(102, 1).Value = 0
$ws.Cells.Item(102, 2).Value = 1
$ws.Cells.Item(102, 3).Value = 0.00655196
$ws.Cells.Item(102, 4).Value = -1
$ws.Cells.Item(102, 5).Value = -1
$ws.Cells.Item(102, 6).Value = 9999
$ws.Cells.Item(102, 7).Value = 1
$ws.Cells.Item(102, 8).Value = -1
$ws.Cells.Item(102, 9).Value = 0
$ws.Cells.Item(102, 10).Value = -1
$ws.Cells.Item(102, 11).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = "no-right"
$ws.Cells.Item(103, 1).Value = 1
$ws.Cells.Item(103, 2).Value = 1
$ws.Cells.Item(103, 3).Value = -0.012027858
$ws.Cells.Item(103, 4).Value = -1
$ws.Cells.Item(103, 5).Value = -1
$ws.Cells.Item(103, 6).Value = 9999
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = -1
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = -1
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = "no-right"
$ws.Cells.Item(104, 1).Value = 1
$ws.Cells.Item(104, 2).Value = 1
$ws.Cells.Item(104, 3).Value = -0.017014052
$ws.Cells.Item(104, 4).Value = -1
$ws.Cells.Item(104, 5).Value = -1
$ws.Cells.Item(104, 6).Value = 9999
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = -1
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = -1
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 13).Value = "no-right"
$ws.Cells.Item(105, 1).Value = 0
$ws.Cells.Item(105, 2).Value = 1
$ws.Cells.Item(105, 3).Value = -0.173381215
$ws.Cells.Item(105, 4).Value = -1
$ws.Cells.Item(105, 5).Value = -1
$ws.Cells.Item(105, 6).Value = 9999
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = -1
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = -1
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = "no-right"
$ws.Cells.Item(106, 1).Value = 0
$ws.Cells.Item(106, 2).Value = 0
$ws.Cells.Item(106, 3).Value = -0.027724113
$ws.Cells.Item(106, 4).Value = -1
$ws.Cells.Item(106, 5).Value = -1
$ws.Cells.Item(106, 6).Value = 9999
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = -1
$ws.Cells.Item(106, 9).Value = 1
$ws.Cells.Item(106, 10).Value = -1
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 13).Value = "no-right"
$ws.Cells.Item(107, 1).Value = 0
$ws.Cells.Item(107, 2).Value = 1
$ws.Cells.Item(107, 3).Value = -0.124972303
$ws.Cells.Item(107, 4).Value = -1
$ws.Cells.Item(107, 5).Value = -1
$ws.Cells.Item(107, 6).Value = 9999
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = -1
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = -1
$ws.Cells.Item(107, 11).Value = 0
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = "no-right"
$ws.Cells.Item(108, 1).Value = 0
$ws.Cells.Item(108, 2).Value = 0
$ws.Cells.Item(108, 3).Value = -0.003798323
$ws.Cells.Item(108, 4).Value = -1
$ws.Cells.Item(108, 5).Value = -1
$ws.Cells.Item(108, 6).Value = 9999
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = -1
$ws.Cells.Item(108, 9).Value = 1
$ws.Cells.Item(108, 10).Value = -1
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 13).Value = "no-right"
$ws.Cells.Item(109, 1).Value = 0
$ws.Cells.Item(109, 2).Value = 1
$ws.Cells.Item(109, 3).Value = -0.015624499
$ws.Cells.Item(109, 4).Value = -1
$ws.Cells.Item(109, 5).Value = -1
$ws.Cells.Item(109, 6).Value = 9999
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = -1
$ws.Cells.Item(109, 9).Value = 0
$ws.Cells.Item(109, 10).Value = -1
$ws.Cells.Item(109, 11).Value = 0
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = "no-right"
$ws.Cells.Item(110, 1).Value = 0
$ws.Cells.Item(110, 2).Value = 0
$ws.Cells.Item(110, 3).Value = 0.033601688
$ws.Cells.Item(110, 4).Value = -1
$ws.Cells.Item(110, 5).Value = -1
$ws.Cells.Item(110, 6).Value = 9999
$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 8).Value = -1
$ws.Cells.Item(110, 9).Value = 1
$ws.Cells.Item(110, 10).Value = -1
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = "no-right"
$ws.Cells.Item(111, 1).Value = 1
$ws.Cells.Item(111, 2).Value = 1
$ws.Cells.Item(111, 3).Value = 0.009619296
$ws.Cells.Item(111, 4).Value = -1
$ws.Cells.Item(111, 5).Value = -1
$ws.Cells.Item(111, 6).Value = 9999
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = -1
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = -1
$ws.Cells.Item(111, 11).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 13).Value = "no-right"
$ws.Cells.Item(112, 1).Value = 0
$ws.Cells.Item(112, 2).Value = 1
$ws.Cells.Item(112, 3).Value = -0.119221941
$ws.Cells.Item(112, 4).Value = -1
$ws.Cells.Item(112, 5).Value = -1
$ws.Cells.Item(112, 6).Value = 9999
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = -1
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = -1
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 13).Value = "no-right"
$ws.Cells.Item(113, 1).Value = 0
$ws.Cells.Item(113, 2).Value = 0
$ws.Cells.Item(113, 3).Value = 0.124013328
$ws.Cells.Item(113, 4).Value = -1
$ws.Cells.Item(113, 5).Value = -1
$ws.Cells.Item(113, 6).Value = 9999
$ws.Cells.Item(113, 7).Value = 1
$ws.Cells.Item(113, 8).Value = -1
$ws.Cells.Item(113, 9).Value = 1
$ws.Cells.Item(113, 10).Value = -1
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = "no-right"
$ws.Cells.Item(114, 1).Value = 1
$ws.Cells.Item(114, 2).Value = 1
$ws.Cells.Item(114, 3).Value = 0.090124937
$ws.Cells.Item(114, 4).Value = -1
$ws.Cells.Item(114, 5).Value = -1
$ws.Cells.Item(114, 6).Value = 9999
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = -1
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = -1
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 13).Value = "no-right"
$ws.Cells.Item(115, 1).Value = 1
$ws.Cells.Item(115, 2).Value = 0
$ws.Cells.Item(115, 3).Value = 0.232528449753817
$ws.Cells.Item(115, 4).Value = 0
$ws.Cells.Item(115, 5).Value = 1
$ws.Cells.Item(115, 6).Value = 0.361346782588241
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 1
$ws.Cells.Item(115, 9).Value = 1
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = "right"
$ws.Cells.Item(116, 1).Value = 1
$ws.Cells.Item(116, 2).Value = 0
$ws.Cells.Item(116, 3).Value = 0.222120029249608
$ws.Cells.Item(116, 4).Value = 1
$ws.Cells.Item(116, 5).Value = 0
$ws.Cells.Item(116, 6).Value = -0.277540461531369
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 1
$ws.Cells.Item(116, 10).Value = 1
$ws.Cells.Item(116, 11).Value = 1
$ws.Cells.Item(116, 12).Value = 1
$ws.Cells.Item(116, 13).Value = "left"
$ws.Cells.Item(117, 1).Value = 1
$ws.Cells.Item(117, 2).Value = 0
$ws.Cells.Item(117, 3).Value = 0.161934428316323
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(117, 5).Value = 1
$ws.Cells.Item(117, 6).Value = 0.355157581181154
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 1
$ws.Cells.Item(117, 9).Value = 1
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 13).Value = "right"
$ws.Cells.Item(118, 1).Value = 1
$ws.Cells.Item(118, 2).Value = 0
$ws.Cells.Item(118, 3).Value = 0.204694473401068
$ws.Cells.Item(118, 4).Value = 1
$ws.Cells.Item(118, 5).Value = 0
$ws.Cells.Item(118, 6).Value = -0.27399005592427
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 9).Value = 1
$ws.Cells.Item(118, 10).Value = 1
$ws.Cells.Item(118, 11).Value = 1
$ws.Cells.Item(118, 12).Value = 1
$ws.Cells.Item(118, 13).Value = "left"

# --- Update view/selection to mirror the target sheetView state ---
$ws.Activate()
$ws.Range("A1:M118").Select()
$excel.ActiveWindow.ScrollRow = 2
